# Update "想去人数" (F column) figures on both the "展览" and "全部类型"
# sheets, which mirror the same underlying dataset.
$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 10188
    5  = 36
    6  = 18
    9  = 17
    12 = 435
    15 = 525
    16 = 19
    17 = 272
    19 = 324
    20 = 109
    21 = 80
    22 = 40
    24 = 19
    25 = 75
    26 = 780
    27 = 1349
    29 = 286
    31 = 55
    32 = 149
    33 = 30
    36 = 186
    37 = 223
    38 = 239
    39 = 1569
    43 = 550
    44 = 3171
    45 = 384
    46 = 161
    47 = 824
    49 = 36
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
